$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reorder the "world" columns: column C (world 1) and column F (world 4)
# swap places with column D (world 2) sliding in between, i.e. a 3-cycle:
#   new C = old F
#   new D = old C
#   new F = old D
# Column E (world 3) and G (branching factor) stay in place.
$tempC = $ws.Range("C2:C7").Value()
$tempD = $ws.Range("D2:D7").Value()
$tempF = $ws.Range("F2:F7").Value()

$ws.Range("C2:C7").Value = $tempF
$ws.Range("D2:D7").Value = $tempC
$ws.Range("F2:F7").Value = $tempD

# Update selection to match the author's final selection (column F was the
# last range they touched).
$ws.Range("F2:F7").Select()
